$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3074.994130329068
$ws.Range("D2").Value = 9990.927465433009
$ws.Range("H2").Value = 52.28189156677908
$ws.Range("I2").Value = 9.829484154314454
$ws.Range("J2").Value = 338.9834683803142
$ws.Range("K2").Value = 35709.00890922388
$ws.Range("L2").Value = 73632.72488415125
$ws.Range("M2").Value = 744853.9651349846
$ws.Range("P2").Value = 356159.9599555623
$ws.Range("L4").Value = 73217.65388798401
$ws.Range("M4").Value = 2545166.61841985
$ws.Range("D5").Value = 10065.87483001023
$ws.Range("F5").Value = 693877.7992514279
$ws.Range("K5").Value = 16451.42835168915
$ws.Range("L5").Value = 44614.97432854872
$ws.Range("P5").Value = 2603925.131103361
$ws.Range("L6").Value = 69713.06049263012
$ws.Range("K7").Value = 10711.99999999995
$ws.Range("L7").Value = 8164.16900863609
$ws.Range("M7").Value = 54366.53152584554
$ws.Range("N7").Value = 1.999999999999993
$ws.Range("K8").Value = 8125.999999999992
$ws.Range("M8").Value = 216908.5926211885
$ws.Range("B9").Value = 4564.627890295273
$ws.Range("D9").Value = 16177.06409222068
$ws.Range("E9").Value = 737.7220073590736
$ws.Range("K9").Value = 4886.999999999982
$ws.Range("L9").Value = 9316.235897955692
$ws.Range("O9").Value = 1143.729881414369
$ws.Range("K10").Value = 2243.999999999904
$ws.Range("L10").Value = 13298.13054705717
$ws.Range("M10").Value = 184575.0998610461
$ws.Range("N10").Value = 5.999999999999972
$ws.Range("D12").Value = 3746.599999999997
$ws.Range("E12").Value = 44.99999999999995
$ws.Range("F12").Value = 84821.61661287429
$ws.Range("I12").Value = 9.999999999999948
$ws.Range("J12").Value = 220.5462165225156
$ws.Range("M12").Value = 633508.489027101
$ws.Range("P12").Value = 286032.9961373269
$ws.Range("E15").Value = 130.9999999999998
$ws.Range("F15").Value = 153719.0900088869
$ws.Range("I15").Value = 2.999999999999976
$ws.Range("B17").Value = 1056.473753859249
$ws.Range("D17").Value = 2381.065679515827
$ws.Range("E17").Value = 113.4249949517034
$ws.Range("F17").Value = 81501.24624526063
$ws.Range("H17").Value = 15.9825267133037
$ws.Range("I17").Value = 2.573634099879445
$ws.Range("J17").Value = 109.9937054332045
$ws.Range("K17").Value = 44335.91368807349
$ws.Range("L17").Value = 2003.262820579908
$ws.Range("M17").Value = 1077875.228372509
$ws.Range("N17").Value = 11.57791462094847
$ws.Range("O17").Value = 315.8466231128759
$ws.Range("P17").Value = 1122077.832194185
$ws.Range("E18").Value = 39.69874823309618
$ws.Range("F18").Value = 377439.3726886988
$ws.Range("J18").Value = 273.7681169627428
$ws.Range("L18").Value = 30552.937259517
$ws.Range("B19").Value = 2528.038589716432
$ws.Range("D19").Value = 11096.21063674203
$ws.Range("E19").Value = 333.1859226706291
$ws.Range("F19").Value = 480263.6654099055
$ws.Range("K19").Value = 253392.3988960179
$ws.Range("L19").Value = 30747.73245897568
$ws.Range("M19").Value = 3041350.614883801
$ws.Range("P19").Value = 482329.4836359872
$ws.Range("M20").Value = 2846566.710468096
$ws.Range("J22").Value = 201.0999999999998
$ws.Range("K22").Value = 14598.99999999995
$ws.Range("L22").Value = 23998.86785081873
$ws.Range("M22").Value = 286065.1040134371
$ws.Range("N22").Value = 10.99999999999997
$ws.Range("D23").Value = 16196.06113338833
$ws.Range("E23").Value = 593.7048669751483
$ws.Range("F23").Value = 198628.8662816901
$ws.Range("J23").Value = 169.2999999999999
$ws.Range("K23").Value = 19625.99999999998
$ws.Range("L23").Value = 14867.20570471301
$ws.Range("M23").Value = 527779.8131531713
$ws.Range("N23").Value = 11.99999999999997
$ws.Range("J24").Value = 188.4999999999993
$ws.Range("K24").Value = 1592.999999999883
$ws.Range("M24").Value = 385565.2847686771
$ws.Range("N24").Value = 7.999999999999967
$ws.Range("J25").Value = 233.2999999999996
$ws.Range("K25").Value = 383.99999999993
$ws.Range("M25").Value = 397106.5589331404
$ws.Range("N25").Value = 0.9999999999999942
$ws.Range("L27").Value = 12579.47779916077
$ws.Range("O27").Value = 801.1041567238514
$ws.Range("P27").Value = 1155335.285561846
$ws.Range("J29").Value = 257.8755753993303
$ws.Range("K29").Value = 2327.949027416018
$ws.Range("M29").Value = 352889.0944825563
$ws.Range("N29").Value = 11.99999999999995
$ws.Range("E31").Value = 104.0508328856842
$ws.Range("K31").Value = 19594.91258346423
$ws.Range("M31").Value = 708166.2951166953
$ws.Range("N31").Value = 5.999999999999963
